# "Pushing the stack implementations"
#
# Marks a handful of additional "Basic Problems" rows as Done (column C)
# on the "Basic Problems" sheet, and updates the saved cursor/scroll
# position on the "Basic Problems" and "Intermediate Problems" sheets.

$wb = $excel.ActiveWorkbook

# --- Basic Problems sheet: mark rows 17,18,19,21,23,26,27,28,29 as Done ---
$basic = $wb.Worksheets.Item("Basic Problems")

$doneRows = @(17, 18, 19, 21, 23, 26, 27, 28, 29)
foreach ($r in $doneRows) {
    $basic.Range("C$r").Value = "Done"
}

# --- Update view state (active sheet, scroll position, selection) ---

# Basic Problems ends up the active sheet, scrolled so row 13 is at the
# top, with E28 selected.
$basic.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$basic.Range("E28").Select()

# Intermediate Problems keeps its B136 selection, just scrolled up so
# row 30 (instead of row 114) is at the top.
$intermediate = $wb.Worksheets.Item("Intermediate Problems")
$intermediate.Activate()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$intermediate.Range("B136").Select()

# Restore Basic Problems as the selected/active tab.
$basic.Activate()
